# Daily attendance processing - 2025-12-03 12:45:24
# Reverses the order of the comma-separated "Recorded By" entries (column G)
# for every row whose value begins with "System, " (e.g. "System, x@y.com"
# becomes "x@y.com, System"). Rows that don't start with "System, " (single
# value, or already starting with something else) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System, ")) {
        $parts = $val.Split(",")
        $n = $parts.Length
        $reversed = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i].Trim()
        }
        $newVal = [string]::Join(", ", $reversed)
        $cell.Value = $newVal
    }
}
